# Adds a new "intervention_type" column (K) to the clinical-trials list.
# Header in K1 gets the same style as the other header cells (A1:J1),
# and K2:K90 get the intervention type for each trial row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, matching the style used by the rest of row 1
# (bold font, thin border all around, centered/top aligned).
$ws.Range("K1").Value = "intervention_type"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("K1").VerticalAlignment = -4160    # xlTop
$ws.Range("K1").Borders.LineStyle = 1

# Intervention type per trial row (rows 2-90).
$interventionTypes = @(
    "PROCEDURE","DRUG","DRUG","DRUG","DEVICE","BIOLOGICAL","DRUG","OTHER","DRUG","PROCEDURE",
    "PROCEDURE","PROCEDURE","DRUG","DRUG","DEVICE","DRUG","RADIATION","DRUG","DRUG","DEVICE",
    "DEVICE","DEVICE","DRUG","OTHER","OTHER","DRUG","OTHER","DRUG","DEVICE","OTHER",
    "OTHER","DRUG","OTHER","DEVICE","DEVICE","OTHER","DEVICE","DRUG","DIAGNOSTIC_TEST","OTHER",
    "DEVICE","PROCEDURE","DEVICE","PROCEDURE","DEVICE","OTHER","DEVICE","OTHER","OTHER","DEVICE",
    "OTHER","DEVICE","OTHER","DEVICE","DIAGNOSTIC_TEST","DIAGNOSTIC_TEST","DRUG","DRUG","PROCEDURE","DIAGNOSTIC_TEST",
    "OTHER","DEVICE","PROCEDURE","DIAGNOSTIC_TEST","DEVICE","DEVICE","OTHER","DEVICE","GENETIC","OTHER",
    "DRUG","OTHER","DRUG","DEVICE","DRUG","DRUG","OTHER","PROCEDURE","OTHER","DRUG",
    "PROCEDURE","DRUG","OTHER","OTHER","PROCEDURE","DIAGNOSTIC_TEST","DIAGNOSTIC_TEST","GENETIC","DRUG"
)

$startRow = 2
for ($i = 0; $i -lt $interventionTypes.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 11).Value = $interventionTypes[$i]
}
